$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusterNames = @{
    20 = "ECs"
    21 = "FAPs"
    22 = "Neutro"
    23 = "sCs"
}

# Row 2
$ws.Cells.Item(2,1).Value = $clusterNames[20]
$ws.Cells.Item(2,2).Value = "Bmp4"
$ws.Cells.Item(2,3).Value = "Bmpr1b"
$ws.Cells.Item(2,4).Value = $clusterNames[20]
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 9.705785499999999
$ws.Cells.Item(2,8).Value = 19.411571
$ws.Cells.Item(2,9).Value = 0.3626960398572868
$ws.Cells.Item(2,10).Value = 0.3021399662106912
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.5
$ws.Cells.Item(2,13).Value = 0.050199
$ws.Cells.Item(2,14).Value = 0.100398
$ws.Cells.Item(2,15).Value = 0.03278033019113316
$ws.Cells.Item(2,16).Value = 0.02558313640211376
$ws.Cells.Item(2,17).Value = 0.4872207263145
$ws.Cells.Item(2,18).Value = 1.948882905258
$ws.Cells.Item(2,19).Value = 0.01188929594553826
$ws.Cells.Item(2,20).Value = 0.007729687968098154

# Row 3
$ws.Cells.Item(3,1).Value = $clusterNames[20]
$ws.Cells.Item(3,2).Value = "Bmp4"
$ws.Cells.Item(3,3).Value = "Bmpr1b"
$ws.Cells.Item(3,4).Value = $clusterNames[21]
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 9.705785499999999
$ws.Cells.Item(3,8).Value = 19.411571
$ws.Cells.Item(3,9).Value = 0.3626960398572868
$ws.Cells.Item(3,10).Value = 0.3021399662106912
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.7708836666666666
$ws.Cells.Item(3,14).Value = 2.312651
$ws.Cells.Item(3,15).Value = 0.5033929188287568
$ws.Cells.Item(3,16).Value = 0.5893032329676366
$ws.Cells.Item(3,17).Value = 7.482031514120166
$ws.Cells.Item(3,18).Value = 44.89218908472099
$ws.Cells.Item(3,19).Value = 0.1825786181513907
$ws.Cells.Item(3,20).Value = 0.1780520588966928

# Row 4
$ws.Cells.Item(4,1).Value = $clusterNames[20]
$ws.Cells.Item(4,2).Value = "Bmp4"
$ws.Cells.Item(4,3).Value = "Bmpr1b"
$ws.Cells.Item(4,4).Value = $clusterNames[22]
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 9.705785499999999
$ws.Cells.Item(4,8).Value = 19.411571
$ws.Cells.Item(4,9).Value = 0.3626960398572868
$ws.Cells.Item(4,10).Value = 0.3021399662106912
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.09074700000000001
$ws.Cells.Item(4,14).Value = 0.272241
$ws.Cells.Item(4,15).Value = 0.05925848371192178
$ws.Cells.Item(4,16).Value = 0.06937168705798773
$ws.Cells.Item(4,17).Value = 0.8807709167685001
$ws.Cells.Item(4,18).Value = 5.284625500611
$ws.Cells.Item(4,19).Value = 0.02149281737026157
$ws.Cells.Item(4,20).Value = 0.02095995918367905

# Row 5
$ws.Cells.Item(5,1).Value = $clusterNames[20]
$ws.Cells.Item(5,2).Value = "Bmp4"
$ws.Cells.Item(5,3).Value = "Bmpr1b"
$ws.Cells.Item(5,4).Value = $clusterNames[23]
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 9.705785499999999
$ws.Cells.Item(5,8).Value = 19.411571
$ws.Cells.Item(5,9).Value = 0.3626960398572868
$ws.Cells.Item(5,10).Value = 0.3021399662106912
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.6195459999999999
$ws.Cells.Item(5,14).Value = 1.239092
$ws.Cells.Item(5,15).Value = 0.4045682672681883
$ws.Cells.Item(5,16).Value = 0.3157419435722618
$ws.Cells.Item(5,17).Value = 6.013180583382999
$ws.Cells.Item(5,18).Value = 24.052722333532
$ws.Cells.Item(5,19).Value = 0.1467353083900963
$ws.Cells.Item(5,20).Value = 0.09539826016222112

# Row 6
$ws.Cells.Item(6,1).Value = $clusterNames[21]
$ws.Cells.Item(6,2).Value = "Bmp4"
$ws.Cells.Item(6,3).Value = "Bmpr1b"
$ws.Cells.Item(6,4).Value = $clusterNames[20]
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 4.987206333333334
$ws.Cells.Item(6,8).Value = 14.961619
$ws.Cells.Item(6,9).Value = 0.1863671917178862
$ws.Cells.Item(6,10).Value = 0.2328767238425594
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.5
$ws.Cells.Item(6,13).Value = 0.050199
$ws.Cells.Item(6,14).Value = 0.100398
$ws.Cells.Item(6,15).Value = 0.03278033019113316
$ws.Cells.Item(6,16).Value = 0.02558313640211376
$ws.Cells.Item(6,17).Value = 0.250352770727
$ws.Cells.Item(6,18).Value = 1.502116624362
$ws.Cells.Item(6,19).Value = 0.006109178081306527
$ws.Cells.Item(6,20).Value = 0.005957716990941575

# Row 7
$ws.Cells.Item(7,1).Value = $clusterNames[21]
$ws.Cells.Item(7,2).Value = "Bmp4"
$ws.Cells.Item(7,3).Value = "Bmpr1b"
$ws.Cells.Item(7,4).Value = $clusterNames[21]
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 4.987206333333334
$ws.Cells.Item(7,8).Value = 14.961619
$ws.Cells.Item(7,9).Value = 0.1863671917178862
$ws.Cells.Item(7,10).Value = 0.2328767238425594
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.7708836666666666
$ws.Cells.Item(7,14).Value = 2.312651
$ws.Cells.Item(7,15).Value = 0.5033929188287568
$ws.Cells.Item(7,16).Value = 0.5893032329676366
$ws.Cells.Item(7,17).Value = 3.844555904663222
$ws.Cells.Item(7,18).Value = 34.601003141969
$ws.Cells.Item(7,19).Value = 0.09381592461278525
$ws.Cells.Item(7,20).Value = 0.1372350062433318

# Row 8
$ws.Cells.Item(8,1).Value = $clusterNames[21]
$ws.Cells.Item(8,2).Value = "Bmp4"
$ws.Cells.Item(8,3).Value = "Bmpr1b"
$ws.Cells.Item(8,4).Value = $clusterNames[22]
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.987206333333334
$ws.Cells.Item(8,8).Value = 14.961619
$ws.Cells.Item(8,9).Value = 0.1863671917178862
$ws.Cells.Item(8,10).Value = 0.2328767238425594
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.09074700000000001
$ws.Cells.Item(8,14).Value = 0.272241
$ws.Cells.Item(8,15).Value = 0.05925848371192178
$ws.Cells.Item(8,16).Value = 0.06937168705798773
$ws.Cells.Item(8,17).Value = 0.4525740131310001
$ws.Cells.Item(8,18).Value = 4.073166118179
$ws.Cells.Item(8,19).Value = 0.01104383719485096
$ws.Cells.Item(8,20).Value = 0.01615505120949546

# Row 9
$ws.Cells.Item(9,1).Value = $clusterNames[21]
$ws.Cells.Item(9,2).Value = "Bmp4"
$ws.Cells.Item(9,3).Value = "Bmpr1b"
$ws.Cells.Item(9,4).Value = $clusterNames[23]
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.987206333333334
$ws.Cells.Item(9,8).Value = 14.961619
$ws.Cells.Item(9,9).Value = 0.1863671917178862
$ws.Cells.Item(9,10).Value = 0.2328767238425594
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.6195459999999999
$ws.Cells.Item(9,14).Value = 1.239092
$ws.Cells.Item(9,15).Value = 0.4045682672681883
$ws.Cells.Item(9,16).Value = 0.3157419435722618
$ws.Cells.Item(9,17).Value = 3.089803734991333
$ws.Cells.Item(9,18).Value = 18.538822409948
$ws.Cells.Item(9,19).Value = 0.07539825182894347
$ws.Cells.Item(9,20).Value = 0.07352894939879059

# Row 10
$ws.Cells.Item(10,1).Value = $clusterNames[22]
$ws.Cells.Item(10,2).Value = "Bmp4"
$ws.Cells.Item(10,3).Value = "Bmpr1b"
$ws.Cells.Item(10,4).Value = $clusterNames[20]
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 5.739525
$ws.Cells.Item(10,8).Value = 17.218575
$ws.Cells.Item(10,9).Value = 0.2144806299461176
$ws.Cells.Item(10,10).Value = 0.2680061118544322
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.5
$ws.Cells.Item(10,13).Value = 0.050199
$ws.Cells.Item(10,14).Value = 0.100398
$ws.Cells.Item(10,15).Value = 0.03278033019113316
$ws.Cells.Item(10,16).Value = 0.02558313640211376
$ws.Cells.Item(10,17).Value = 0.288118415475
$ws.Cells.Item(10,18).Value = 1.72871049285
$ws.Cells.Item(10,19).Value = 0.007030745869235979
$ws.Cells.Item(10,20).Value = 0.006856436916172095

# Row 11
$ws.Cells.Item(11,1).Value = $clusterNames[22]
$ws.Cells.Item(11,2).Value = "Bmp4"
$ws.Cells.Item(11,3).Value = "Bmpr1b"
$ws.Cells.Item(11,4).Value = $clusterNames[21]
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 5.739525
$ws.Cells.Item(11,8).Value = 17.218575
$ws.Cells.Item(11,9).Value = 0.2144806299461176
$ws.Cells.Item(11,10).Value = 0.2680061118544322
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.7708836666666666
$ws.Cells.Item(11,14).Value = 2.312651
$ws.Cells.Item(11,15).Value = 0.5033929188287568
$ws.Cells.Item(11,16).Value = 0.5893032329676366
$ws.Cells.Item(11,17).Value = 4.424506076925
$ws.Cells.Item(11,18).Value = 39.820554692325
$ws.Cells.Item(11,19).Value = 0.1079680303408066
$ws.Cells.Item(11,20).Value = 0.1579368681709029

# Row 12
$ws.Cells.Item(12,1).Value = $clusterNames[22]
$ws.Cells.Item(12,2).Value = "Bmp4"
$ws.Cells.Item(12,3).Value = "Bmpr1b"
$ws.Cells.Item(12,4).Value = $clusterNames[22]
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 5.739525
$ws.Cells.Item(12,8).Value = 17.218575
$ws.Cells.Item(12,9).Value = 0.2144806299461176
$ws.Cells.Item(12,10).Value = 0.2680061118544322
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.3333333333333333
$ws.Cells.Item(12,13).Value = 0.09074700000000001
$ws.Cells.Item(12,14).Value = 0.272241
$ws.Cells.Item(12,15).Value = 0.05925848371192178
$ws.Cells.Item(12,16).Value = 0.06937168705798773
$ws.Cells.Item(12,17).Value = 0.5208446751750001
$ws.Cells.Item(12,18).Value = 4.687602076575001
$ws.Cells.Item(12,19).Value = 0.01270979691618473
$ws.Cells.Item(12,20).Value = 0.01859203612119372

# Row 13
$ws.Cells.Item(13,1).Value = $clusterNames[22]
$ws.Cells.Item(13,2).Value = "Bmp4"
$ws.Cells.Item(13,3).Value = "Bmpr1b"
$ws.Cells.Item(13,4).Value = $clusterNames[23]
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 5.739525
$ws.Cells.Item(13,8).Value = 17.218575
$ws.Cells.Item(13,9).Value = 0.2144806299461176
$ws.Cells.Item(13,10).Value = 0.2680061118544322
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.6195459999999999
$ws.Cells.Item(13,14).Value = 1.239092
$ws.Cells.Item(13,15).Value = 0.4045682672681883
$ws.Cells.Item(13,16).Value = 0.3157419435722618
$ws.Cells.Item(13,17).Value = 3.55589975565
$ws.Cells.Item(13,18).Value = 21.3353985339
$ws.Cells.Item(13,19).Value = 0.08677205681989031
$ws.Cells.Item(13,20).Value = 0.0846207706461634

# Row 14
$ws.Cells.Item(14,1).Value = $clusterNames[23]
$ws.Cells.Item(14,2).Value = "Bmp4"
$ws.Cells.Item(14,3).Value = "Bmpr1b"
$ws.Cells.Item(14,4).Value = $clusterNames[20]
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 6.327592
$ws.Cells.Item(14,8).Value = 12.655184
$ws.Cells.Item(14,9).Value = 0.2364561384787094
$ws.Cells.Item(14,10).Value = 0.1969771980923172
$ws.Cells.Item(14,11).Value = 1
$ws.Cells.Item(14,12).Value = 0.5
$ws.Cells.Item(14,13).Value = 0.050199
$ws.Cells.Item(14,14).Value = 0.100398
$ws.Cells.Item(14,15).Value = 0.03278033019113316
$ws.Cells.Item(14,16).Value = 0.02558313640211376
$ws.Cells.Item(14,17).Value = 0.317638790808
$ws.Cells.Item(14,18).Value = 1.270555163232
$ws.Cells.Item(14,19).Value = 0.0077511102950524
$ws.Cells.Item(14,20).Value = 0.005039294526901933

# Row 15
$ws.Cells.Item(15,1).Value = $clusterNames[23]
$ws.Cells.Item(15,2).Value = "Bmp4"
$ws.Cells.Item(15,3).Value = "Bmpr1b"
$ws.Cells.Item(15,4).Value = $clusterNames[21]
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 6.327592
$ws.Cells.Item(15,8).Value = 12.655184
$ws.Cells.Item(15,9).Value = 0.2364561384787094
$ws.Cells.Item(15,10).Value = 0.1969771980923172
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.7708836666666666
$ws.Cells.Item(15,14).Value = 2.312651
$ws.Cells.Item(15,15).Value = 0.5033929188287568
$ws.Cells.Item(15,16).Value = 0.5893032329676366
$ws.Cells.Item(15,17).Value = 4.877837322130667
$ws.Cells.Item(15,18).Value = 29.267023932784
$ws.Cells.Item(15,19).Value = 0.1190303457237742
$ws.Cells.Item(15,20).Value = 0.1160792996567091

# Row 16
$ws.Cells.Item(16,1).Value = $clusterNames[23]
$ws.Cells.Item(16,2).Value = "Bmp4"
$ws.Cells.Item(16,3).Value = "Bmpr1b"
$ws.Cells.Item(16,4).Value = $clusterNames[22]
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 6.327592
$ws.Cells.Item(16,8).Value = 12.655184
$ws.Cells.Item(16,9).Value = 0.2364561384787094
$ws.Cells.Item(16,10).Value = 0.1969771980923172
$ws.Cells.Item(16,11).Value = 1
$ws.Cells.Item(16,12).Value = 0.3333333333333333
$ws.Cells.Item(16,13).Value = 0.09074700000000001
$ws.Cells.Item(16,14).Value = 0.272241
$ws.Cells.Item(16,15).Value = 0.05925848371192178
$ws.Cells.Item(16,16).Value = 0.06937168705798773
$ws.Cells.Item(16,17).Value = 0.5742099912240001
$ws.Cells.Item(16,18).Value = 3.445259947344
$ws.Cells.Item(16,19).Value = 0.01401203223062452
$ws.Cells.Item(16,20).Value = 0.01366464054361949

# Row 17
$ws.Cells.Item(17,1).Value = $clusterNames[23]
$ws.Cells.Item(17,2).Value = "Bmp4"
$ws.Cells.Item(17,3).Value = "Bmpr1b"
$ws.Cells.Item(17,4).Value = $clusterNames[23]
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 6.327592
$ws.Cells.Item(17,8).Value = 12.655184
$ws.Cells.Item(17,9).Value = 0.2364561384787094
$ws.Cells.Item(17,10).Value = 0.1969771980923172
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.6195459999999999
$ws.Cells.Item(17,14).Value = 1.239092
$ws.Cells.Item(17,15).Value = 0.4045682672681883
$ws.Cells.Item(17,16).Value = 0.3157419435722618
$ws.Cells.Item(17,17).Value = 3.920234313232
$ws.Cells.Item(17,18).Value = 15.680937252928
$ws.Cells.Item(17,19).Value = 0.09566265022925823
$ws.Cells.Item(17,20).Value = 0.06219396336508665
